# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

$wsRushing = $wb.Worksheets.Item("Rushing")
$wsReceiving = $wb.Worksheets.Item("Receiving")

# ---- Rushing sheet updates ----
# Row 2: J.Allen
$wsRushing.Range("C2").Value = 29
$wsRushing.Range("D2").Value = 31
$wsRushing.Range("E2").Value = 32
$wsRushing.Range("F2").Value = 24

# Row 3: D.Singletary
$wsRushing.Range("C3").Value = 79
$wsRushing.Range("D3").Value = 57
$wsRushing.Range("F3").Value = 24

# Row 4: Z.Moss
$wsRushing.Range("C4").Value = 43
$wsRushing.Range("D4").Value = 30
$wsRushing.Range("F4").Value = 26

# Row 8: I.McKenzie
$wsRushing.Range("C8").Value = 3
$wsRushing.Range("F8").Value = 3

# ---- Receiving sheet updates ----
# Row 2: D.Singletary
$wsReceiving.Range("C2").Value = 44
$wsReceiving.Range("D2").Value = 36

# Row 6: S.Diggs
$wsReceiving.Range("C6").Value = 99
$wsReceiving.Range("D6").Value = 71
$wsReceiving.Range("E6").Value = 30
$wsReceiving.Range("F6").Value = 11

# Row 7: E.Sanders
$wsReceiving.Range("C7").Value = 45
$wsReceiving.Range("D7").Value = 31

# Row 8: C.Beasley
$wsReceiving.Range("C8").Value = 93
$wsReceiving.Range("D8").Value = 71
$wsReceiving.Range("E8").Value = 9
$wsReceiving.Range("F8").Value = 5
$wsReceiving.Range("G8").Value = 13
$wsReceiving.Range("H8").Value = 9

# Row 9: G.Davis
$wsReceiving.Range("C9").Value = 26
$wsReceiving.Range("D9").Value = 16
$wsReceiving.Range("E9").Value = 18
$wsReceiving.Range("F9").Value = 12

# Update selection to match the author's final cursor position
$wsReceiving.Range("E2").Select()
